# Automatische test-sync: 2025-08-19 20:59:50
# Append a new log entry to the "Logs" sheet and refresh the summary
# count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row appended as row 20 on the Logs sheet.
$newRow = 20

$logs.Cells.Item($newRow, 1).Value = "Vraag over product"
$logs.Cells.Item($newRow, 2).Value = "documentatie@testbedrijf123.nl"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 6).Value = "2025-08-19 20:59:11"
$logs.Cells.Item($newRow, 7).Value = "Nee"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Refresh the Dashboard count for the matching category.
$dashboard.Cells.Item(2, 2).Value = 19

# Extend the conditional-formatting ranges so they keep covering the
# full data range now that a row was appended (D/G/H/I/J 2:19 -> 2:20).
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "19")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "20")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
